# The workbook gains one new data row. A new row is inserted at row 59,
# pushing the existing rows 59-181 down to 60-182 (dimension grows from
# A1:R181 to A1:R182). The new row 59 is populated with a new weekly
# price observation for Berenjena / Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 59, shifting rows 59:181 down to 60:182.
$ws.Rows(59).Insert()

# Populate the newly inserted row 59 with the new record.
$ws.Range("A59").Value = 3
$ws.Range("B59").Value = "Femacal de La Calera"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44519
$ws.Range("E59").Value = 5
$ws.Range("F59").Value = 100112001
$ws.Range("G59").Value = "Berenjena"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 125
$ws.Range("K59").Value = 7000
$ws.Range("L59").Value = 7500
$ws.Range("M59").Value = 7240
$ws.Range("N59").Value = "`$/caja 60 unidades"
$ws.Range("O59").Value = "Región de Arica y Parinacota"
$ws.Range("P59").Value = 121
$ws.Range("Q59").Value = 60
$ws.Range("R59").Value = "Hortaliza"
